# QuickStart - 6 - Property & Event Binding
# 1) Drop the old "_GoBack" bookmark that sat after "...property to false".
# 2) Append "POINTS TO NOTE" section with a 4-item numbered list, the
#    "_GoBack" bookmark now trailing the first bullet's "click event" text.
# 3) Along the way Word (re)materialises the "List Paragraph" style, a
#    fresh numbered-list definition (numId 4) and the "Balloon Text" /
#    "Balloon Text Char" styles that a paste-from-another-doc session
#    tends to drag in.

$d = $word.ActiveDocument

# --- 1. remove the stale _GoBack bookmark -----------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. materialise the styles this edit references -------------------
$balloonText = $d.Styles.Add("Balloon Text", 1)
$balloonText.BaseStyle = $d.Styles("Normal")
$balloonText.LinkStyle = "BalloonTextChar"
$balloonText.Priority = 99
$balloonText.UnhideWhenUsed = $true
$balloonText.ParagraphFormat.SpaceAfter = 0
$balloonText.ParagraphFormat.LineSpacingRule = 0
$balloonText.Font.NameAscii = "Segoe UI"
$balloonText.Font.NameOther = "Segoe UI"
$balloonText.Font.NameBi = "Segoe UI"
$balloonText.Font.Size = 9
$balloonText.Font.SizeBi = 9

$balloonTextChar = $d.Styles.Add("Balloon Text Char", 2)
$balloonTextChar.BaseStyle = $d.Styles("DefaultParagraphFont")
$balloonTextChar.LinkStyle = "BalloonText"
$balloonTextChar.Priority = 99
$balloonTextChar.Font.NameAscii = "Segoe UI"
$balloonTextChar.Font.NameFarEast = "Calibri"
$balloonTextChar.Font.NameOther = "Segoe UI"
$balloonTextChar.Font.NameBi = "Segoe UI"
$balloonTextChar.Font.Color = 3355443
$balloonTextChar.Font.Size = 9
$balloonTextChar.Font.SizeBi = 9

$listParagraph = $d.Styles.Add("List Paragraph", 1)
$listParagraph.BaseStyle = $d.Styles("Normal")
$listParagraph.Priority = 34
$listParagraph.QuickStyle = $true
$listParagraph.ParagraphFormat.LeftIndent = 36

# --- 3. materialise the "1) 2) 3) ..." numbered list (numId 4) --------
# Borrow the very last (empty) paragraph to mint the numbering
# definition, then strip the numbering back off so that paragraph is
# left exactly as it was.
$tail = $d.Paragraphs.Last
$tail.Range.ListFormat.ApplyNumberDefault()
$tail2 = $d.Paragraphs.Last
$tail2.Range.ListFormat.RemoveNumbers()

# --- 4. append the new "POINTS TO NOTE" block --------------------------
$insertAt = $d.Paragraphs.Last.Range
$newContentXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>POINTS TO NOTE</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>click</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> event</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>hidden</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>DOM property</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t>class doesn&#8217;t have export keyword</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">while accessing the class variable remember to use </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>this</w:t>
  </w:r>
</w:p>
"@
$insertAt.InsertXML($newContentXml)

Write-Host "Applied property & event binding POINTS TO NOTE edit."
